# 20150723 - Code cleanup, fomatting cleanup, add SPI1 and SPI2 pin definitions.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pins")

# Add SD Card CS pin definition (row 12, PE10)
$ws.Range("J12").Value = "SD Card CS"
$ws.Range("K12").Value = "SD Card CS"

# Add SPI1 and SPI2 CS pin definitions (row 5, PA3)
$ws.Range("B5").Value = "SP12 CS"
$ws.Range("C5").Value = "SPI2 CS"

# Duplicate AQ32 Definition into TauLabs Definition column for SPI1 pins (rows 7-9)
$ws.Range("C7").Value = $ws.Range("B7").Value2
$ws.Range("C8").Value = $ws.Range("B8").Value2
$ws.Range("C9").Value = $ws.Range("B9").Value2

# Duplicate AQ32 Definition into TauLabs Definition column for SPI2 pins (rows 32-34)
$ws.Range("C32").Value = $ws.Range("B32").Value2
$ws.Range("C33").Value = $ws.Range("B33").Value2
$ws.Range("C34").Value = $ws.Range("B34").Value2

# Update the view: move the active selection
$ws.Range("G26").Select()
